$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update text content (DoNotThrow branding / copy edits) ---
$ws.Range("D9").Value  = "Homepage of the website with information on Donothrow"
$ws.Range("E7").Value  = "Presentation of Donothrow"
$ws.Range("F9").Value  = "Presentation of Donothrow, link to the website"
$ws.Range("F11").Value = "Presentation of Donothrow, link to the website"
$ws.Range("J9").Value  = "Inserts address and type of food needed"
$ws.Range("J11").Value = "Inserts address and type of food needed"
$ws.Range("M9").Value  = "Checks some reasons of departure in a list"
$ws.Range("M11").Value = "Checks some reasons of departure in a list"

# M14 loses its content entirely (the "insertion of a new reason of departure" DB step was removed)
$ws.Range("M14").Value = ""

# --- Row 9 shrank after the copy edits (autofit-style height change) ---
$ws.Rows(9).RowHeight = 112.8

# --- Selection moved to E7 ---
$ws.Range("E7").Select()
